$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.885.87'
$ws.Range("E2").Value = '  +1.97%  '
$ws.Range("D3").Value = '1.881.83'
$ws.Range("E3").Value = '  +1.65%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '333.02'
$ws.Range("E5").Value = '  +3.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  +0.11%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4732'
$ws.Range("E7").Value = '  +5.79%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3973'
$ws.Range("E8").Value = '  +4.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.25'
$ws.Range("E9").Value = '  +0.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08056'
$ws.Range("E10").Value = '  +2.51%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.029'
$ws.Range("E11").Value = '  +1.74%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.89'
$ws.Range("E12").Value = '  +2.65%  '
$ws.Range("D13").Value = '1.911.83'
$ws.Range("E13").Value = '  +3.42%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.963'
$ws.Range("E14").Value = '  +1.97%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.198'
$ws.Range("E15").Value = '  +1.34%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.003'
$ws.Range("E16").Value = '  -0.04%  '
$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '87.29'
$ws.Range("E17").Value = '  +1.71%  '
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001048'
$ws.Range("E18").Value = '  +1.57%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06618'
$ws.Range("E19").Value = '  +1.80%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.31'
$ws.Range("E20").Value = '  +2.24%  '
$ws.Range("E21").Value = '  +0.10%  '
$ws.Range("D22").Value = '27.936.57'
$ws.Range("E22").Value = '  +2.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.519'
$ws.Range("E23").Value = '  +1.13%  '
$ws.Range("E24").Value = '  +2.87%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.313'
$ws.Range("E25").Value = '  +2.74%  '
$ws.Range("D26").Value = '2.120.03'
$ws.Range("E26").Value = '  +2.50%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '157.83'
$ws.Range("E27").Value = '  +4.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.23'
$ws.Range("E28").Value = '  +4.48%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.109'
$ws.Range("E29").Value = '  +2.60%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.625'
$ws.Range("E30").Value = '  +1.95%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '122.65'
$ws.Range("E31").Value = '  +2.63%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9863'
$ws.Range("E32").Value = '  +5.72%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09576'
$ws.Range("E33").Value = '  +2.82%  '
$ws.Range("E34").Value = '  +0.10%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.626'
$ws.Range("E35").Value = '  +0.61%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.334'
$ws.Range("E36").Value = '  +1.76%  '
$ws.Range("E37").Value = '  +2.86%  '
$ws.Range("E38").Value = '  +2.17%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.229'
$ws.Range("E39").Value = '  +2.22%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.248'
$ws.Range("E40").Value = '  -0.26%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6032'
$ws.Range("E41").Value = '  +2.40%  '
$ws.Range("E42").Value = '  +0.12%  '
$ws.Range("E43").Value = '  +3.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.36'
$ws.Range("E44").Value = '  +1.11%  '
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5721'
$ws.Range("E45").Value = '  +1.67%  '
$ws.Range("B46").Value = 'WEMIXTOKEN'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.254'
$ws.Range("E46").Value = '  +0.05%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.28'
$ws.Range("E47").Value = '  +0.50%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.416'
$ws.Range("E48").Value = '  +1.69%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.952'
$ws.Range("E49").Value = '  +1.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06839'
$ws.Range("E50").Value = '  -0.27%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '113.94'
$ws.Range("E51").Value = '  +5.29%  '
